$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H19").Value = 1124.3636
$ws.Range("I19").Value = 1740
$ws.Range("J19").Value = 893.5
$ws.Range("K19").Value = 1740
$ws.Range("L19").Value = 893.5
$ws.Range("M19").Value = -1565
$ws.Range("N19").Value = -1243.5
$ws.Range("H129").Value = 2015.1505
$ws.Range("I129").Value = 5771.2104
$ws.Range("J129").Value = 1050.7567
$ws.Range("K129").Value = 17313.6312
$ws.Range("L129").Value = 3152.2701
$ws.Range("M129").Value = -12313.6312
$ws.Range("N129").Value = -13152.2701
$ws.Range("H132").Value = 4314454.5
$ws.Range("I132").Value = 4721365.5
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 14164096.5
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -14161566.5
$ws.Range("N132").Value = -8660

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 8411.223
$ws.Range("I6").Value = 17234
$ws.Range("J6").Value = 3999.8333
$ws.Range("K6").Value = 17234
$ws.Range("L6").Value = 3999.8333
$ws.Range("M6").Value = -17061
$ws.Range("N6").Value = -4345.8333
$ws.Range("H45").Value = 84849.414
$ws.Range("I45").Value = 101528.2
$ws.Range("K45").Value = 101528.2
$ws.Range("M45").Value = -101151.2

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1507.4688
$ws.Range("I99").Value = 1290.8334
$ws.Range("J99").Value = 1557.4615
$ws.Range("K99").Value = 1290.8334
$ws.Range("L99").Value = 1557.4615
$ws.Range("M99").Value = 207.1666
$ws.Range("N99").Value = -4553.461499999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2292.2144
$ws.Range("I58").Value = 1766.6666
$ws.Range("J58").Value = 2686.375
$ws.Range("K58").Value = 1766.6666
$ws.Range("L58").Value = 2686.375
$ws.Range("M58").Value = -1563.6666
$ws.Range("N58").Value = -3092.375
$ws.Range("H105").Value = 1125.8572
$ws.Range("I105").Value = 1041.4375
$ws.Range("J105").Value = 1396
$ws.Range("K105").Value = 1041.4375
$ws.Range("L105").Value = 1396
$ws.Range("M105").Value = 705.5625
$ws.Range("N105").Value = -4890
$ws.Range("H136").Value = 2292.2144
$ws.Range("I136").Value = 1766.6666
$ws.Range("J136").Value = 2686.375
$ws.Range("K136").Value = 5299.9998
$ws.Range("L136").Value = 8059.125
$ws.Range("M136").Value = -2749.9998
$ws.Range("N136").Value = -13159.125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1325.3334
$ws.Range("J5").Value = 1411.1154
$ws.Range("L5").Value = 4233.3462
$ws.Range("N5").Value = -4457.3462
$ws.Range("H68").Value = 1233.25
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1254.4546
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 3763.3638
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -5385.3638
$ws.Range("H71").Value = 1233.25
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1254.4546
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 11290.0914
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -19402.0914
$ws.Range("H97").Value = 1344.8889
$ws.Range("J97").Value = 1326
$ws.Range("L97").Value = 3978
$ws.Range("N97").Value = -4970
$ws.Range("H117").Value = 9044
$ws.Range("I117").Value = 1190
$ws.Range("J117").Value = 10353
$ws.Range("K117").Value = 3570
$ws.Range("L117").Value = 31059
$ws.Range("M117").Value = -128
$ws.Range("N117").Value = -37943
$ws.Range("H132").Value = 1638.5518
$ws.Range("I132").Value = 796.5
$ws.Range("J132").Value = 2232.9412
$ws.Range("K132").Value = 7168.5
$ws.Range("L132").Value = 20096.4708
$ws.Range("M132").Value = -4638.5
$ws.Range("N132").Value = -25156.4708
$ws.Range("H133").Value = 3207.8572
$ws.Range("I133").Value = 1903.75
$ws.Range("K133").Value = 5711.25
$ws.Range("M133").Value = -651.25
$ws.Range("H134").Value = 3937.7646
$ws.Range("I134").Value = 2397.5
$ws.Range("J134").Value = 7634.4
$ws.Range("K134").Value = 7192.5
$ws.Range("L134").Value = 22903.2
$ws.Range("M134").Value = -2122.5
$ws.Range("N134").Value = -33043.2
$ws.Range("H135").Value = 1325.3334
$ws.Range("J135").Value = 1411.1154
$ws.Range("L135").Value = 12700.0386
$ws.Range("N135").Value = -17770.0386
$ws.Range("H136").Value = 1956.6666
$ws.Range("J136").Value = 2144.4443
$ws.Range("L136").Value = 6433.3329
$ws.Range("N136").Value = -16633.3329
$ws.Range("H137").Value = 50010.39
$ws.Range("I137").Value = 51812.3
$ws.Range("J137").Value = 37997.668
$ws.Range("K137").Value = 155436.9
$ws.Range("L137").Value = 113993.004
$ws.Range("M137").Value = -150336.9
$ws.Range("N137").Value = -124193.004

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 217.3
$ws.Range("I2").Value = 155
$ws.Range("K2").Value = 155
$ws.Range("M2").Value = -42
$ws.Range("H26").Value = 9150
$ws.Range("J26").Value = 9150
$ws.Range("L26").Value = 9150
$ws.Range("N26").Value = -9710
$ws.Range("H46").Value = 11407.546
$ws.Range("J46").Value = 12498.111
$ws.Range("L46").Value = 12498.111
$ws.Range("N46").Value = -12810.111
$ws.Range("H50").Value = 9150
$ws.Range("J50").Value = 9150
$ws.Range("L50").Value = 9150
$ws.Range("N50").Value = -10146
$ws.Range("H57").Value = 19500
$ws.Range("J57").Value = 19500
$ws.Range("L57").Value = 19500
$ws.Range("N57").Value = -21140
$ws.Range("H80").Value = 83420880
$ws.Range("J80").Value = 1944.25
$ws.Range("L80").Value = 1944.25
$ws.Range("N80").Value = -3940.25
$ws.Range("H83").Value = 83420880
$ws.Range("J83").Value = 1944.25
$ws.Range("L83").Value = 9721.25
$ws.Range("N83").Value = -19705.25
$ws.Range("H126").Value = 3032.8462
$ws.Range("I126").Value = 3751.5715
$ws.Range("J126").Value = 2194.3333
$ws.Range("K126").Value = 11254.7145
$ws.Range("L126").Value = 6582.999899999999
$ws.Range("M126").Value = -8784.7145
$ws.Range("N126").Value = -11522.9999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5668159.5
$ws.Range("I16").Value = 8400670
$ws.Range("K16").Value = 8400670
$ws.Range("M16").Value = -8400500
$ws.Range("H61").Value = 1631.3
$ws.Range("I61").Value = 1834.1111
$ws.Range("J61").Value = 1465.3636
$ws.Range("K61").Value = 1834.1111
$ws.Range("L61").Value = 1465.3636
$ws.Range("M61").Value = -1632.1111
$ws.Range("N61").Value = -1869.3636
$ws.Range("H68").Value = 3105.2632
$ws.Range("J68").Value = 3422.0715
$ws.Range("L68").Value = 3422.0715
$ws.Range("N68").Value = -4920.0715
$ws.Range("H71").Value = 3105.2632
$ws.Range("J71").Value = 3422.0715
$ws.Range("L71").Value = 17110.3575
$ws.Range("N71").Value = -24598.3575
$ws.Range("H93").Value = 2623.5625
$ws.Range("I93").Value = 2945.889
$ws.Range("J93").Value = 2209.1428
$ws.Range("K93").Value = 2945.889
$ws.Range("L93").Value = 2209.1428
$ws.Range("M93").Value = -1697.889
$ws.Range("N93").Value = -4705.1428
$ws.Range("H113").Value = 1631.3
$ws.Range("I113").Value = 1834.1111
$ws.Range("J113").Value = 1465.3636
$ws.Range("K113").Value = 1834.1111
$ws.Range("L113").Value = 1465.3636
$ws.Range("M113").Value = 335.8888999999999
$ws.Range("N113").Value = -5805.3636
